# Applies the template placeholder rewrite described by the commit
# "Updated templates for documents" to group_template.docx.
#
# Strategy: all of the run-splitting in the diff is purely cosmetic
# (Word re-flowing runs as a human retyped the placeholders); the
# underlying paragraph structure and text stream are unchanged, so a
# sequence of literal Find/Replace operations over $d.Content
# reproduces the target text exactly. A handful of table column widths
# also shifted by 1 dxa (twip); those are applied via the Tables OM
# since they are not representable as a text find/replace.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. "{{practice_kind}} " -> "{{practice.kind.name_dp}} "
Replace-Text "{{practice_kind}} " "{{practice.kind.name_dp}} "

# 2. "{{practice_type}} практики" -> "({{practice.type.name_dp}}) практики"
Replace-Text "{{practice_type}} практики" "({{practice.type.name_dp}}) практики"

# 3. "Курс {{course}}" -> "Курс {{group.name}}"
Replace-Text "Курс {{course}}" "Курс {{group.name}}"

# 4. "Группа {{group}}" -> "Группа {{group.course}}"
Replace-Text "Группа {{group}}" "Группа {{group.course}}"

# 5. "г. Ханты-Мансийск, {{year}}г" -> "г. Ханты-Мансийск, {{practice.year}}г"
Replace-Text "г. Ханты-Мансийск, {{year}}г" "г. Ханты-Мансийск, {{practice.year}}г"

# 6. practice_start / practice_end -> practice.start.date / practice.end.date
Replace-Text "с {{practice_start}} по {{practice_end}}" "с {{practice.start.date}} по {{practice.end.date}}"

# 7. "2. Номер и дата приказа: № {{report_number}}  от {{report_date}}"
#    -> "2. Номер и дата приказа: № {{practice.order}}"
Replace-Text "№ {{report_number}}  от {{report_date}}" "№ {{practice.order}}"

# 8. "3. Вид практики: {{practice_kind_ip}}" -> "3. Вид практики: {{practice.kind.name}}"
Replace-Text "{{practice_kind_ip}}" "{{practice.kind.name}}"

# 9. "4. Тип практики: {{practice_type_ip}}" -> "4. Тип практики: {{practice.type.name}}"
Replace-Text "{{practice_type_ip}}" "{{practice.type.name}}"

# 10. "5. Количество обучающихся прошедших практику: {{students_number_success}}"
#     -> "... {{group.success_students_number}}"
Replace-Text "{{students_number_success}}" "{{group.success_students_number}}"

# 11. "{%tr for student_success in students_success %}"
#     -> "{%tr for student_success in group.success_students %}"
Replace-Text "students_success %}" "group.success_students %}"

# 12. "{{student_success.payable}}" -> "{{student_success.paid}}"
Replace-Text "{{student_success.payable}}" "{{student_success.paid}}"

# 13. "{{student_success.fullname_manager}}" -> "{{student_success.directior_orgainzation}}"
Replace-Text "{{student_success.fullname_manager}}" "{{student_success.directior_orgainzation}}"

# 14. "6. Количество обучающихся не прошедших практику: {{students_number_fail}}"
#     -> "... {{group.failed_students_number}}"
Replace-Text "{{students_number_fail}}" "{{group.failed_students_number}}"

# 15. "{%tr for student_fail in students_fail %}"
#     -> "{%tr for student_fail in group.failed_students %}"
Replace-Text "students_fail %}" "group.failed_students %}"

# 16. "{{student_fail.comment}}" -> "{{student_fail.reason}}"
Replace-Text "{{student_fail.comment}}" "{{student_fail.reason}}"

# 17. "{{recommendation}}" -> "{{practice.recommendation}}"
Replace-Text "{{recommendation}}" "{{practice.recommendation}}"

# 18. "{{usu_manager}}" -> "{{practice.usu_name_short}}"
Replace-Text "{{usu_manager}}" "{{practice.usu_name_short}}"

# 19. "{{opop_manager}}" -> "{{practice.opop_name_short}}"
Replace-Text "{{opop_manager}}" "{{practice.opop_name_short}}"

# Table column width tweaks (students table): 1981/938/1278/1995 dxa
# become 1980/939/1279/1994 dxa. COM widths are expressed in points
# (1 pt = 20 dxa); setting Columns(i).Width rewrites both the tblGrid
# entry and every cell in that column.
$t = $d.Tables.Item(1)
$t.Columns.Item(2).Width = 1980 / 20.0
$t.Columns.Item(3).Width = 939 / 20.0
$t.Columns.Item(6).Width = 1279 / 20.0
$t.Columns.Item(7).Width = 1994 / 20.0

Write-Output "done"
